# Auto-generated edit script applying the 2026-02-21 07:20 meteocat data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Percent-valued cells (column H) must stay literal text ("NN%") rather than
# being auto-parsed into a percentage number by the COM Value setter, so force
# the number format to Text ("@") before writing those specific cells.
$percentCells = @("H2", "H5", "H12", "H14", "H17", "H19", "H20", "H23", "H25", "H27", "H29", "H30", "H31", "H33", "H36", "H40", "H41", "H42", "H46")
foreach ($pc in $percentCells) {
    $ws.Range($pc).NumberFormat = "@"
}

$ws.Range("E2").Value = "2026-02-21 07:18:19"
$ws.Range("H2").Value = "54%"
$ws.Range("E3").Value = "2026-02-21 07:18:21"
$ws.Range("O3").Value = "-1.2 °C"
$ws.Range("E4").Value = "2026-02-21 07:18:24"
$ws.Range("J4").Value = "1029.4 hPa"
$ws.Range("E5").Value = "2026-02-21 07:18:26"
$ws.Range("H5").Value = "44%"
$ws.Range("O5").Value = "0.0 °C"
$ws.Range("E6").Value = "2026-02-21 07:18:28"
$ws.Range("J6").Value = "1029.2 hPa"
$ws.Range("L6").Value = "13.7 km/h - 302º 6:53 TU"
$ws.Range("E7").Value = "2026-02-21 07:18:31"
$ws.Range("J7").Value = "1028.8 hPa"
$ws.Range("O7").Value = "11.0 °C"
$ws.Range("E8").Value = "2026-02-21 07:18:33"
$ws.Range("J8").Value = "1029.0 hPa"
$ws.Range("N8").Value = "6.6 °C 6:59 TU"
$ws.Range("O8").Value = "7.4 °C"
$ws.Range("E9").Value = "2026-02-21 07:18:35"
$ws.Range("O9").Value = "11.3 °C"
$ws.Range("E10").Value = "2026-02-21 07:18:38"
$ws.Range("E11").Value = "2026-02-21 07:18:40"
$ws.Range("O11").Value = "5.3 °C"
$ws.Range("E12").Value = "2026-02-21 07:18:43"
$ws.Range("H12").Value = "58%"
$ws.Range("N12").Value = "7.2 °C 6:39 TU"
$ws.Range("O12").Value = "11.2 °C"
$ws.Range("E13").Value = "2026-02-21 07:18:45"
$ws.Range("J13").Value = "1036.1 hPa"
$ws.Range("O13").Value = "-3.6 °C"
$ws.Range("E14").Value = "2026-02-21 07:18:47"
$ws.Range("H14").Value = "73%"
$ws.Range("L14").Value = "29.5 km/h - 335º 6:46 TU"
$ws.Range("O14").Value = "8.8 °C"
$ws.Range("E15").Value = "2026-02-21 07:18:50"
$ws.Range("O15").Value = "11.0 °C"
$ws.Range("E16").Value = "2026-02-21 07:18:52"
$ws.Range("O16").Value = "0.6 °C"
$ws.Range("E17").Value = "2026-02-21 07:18:54"
$ws.Range("H17").Value = "33%"
$ws.Range("M17").Value = "7.6 °C 6:51 TU"
$ws.Range("O17").Value = "6.6 °C"
$ws.Range("E18").Value = "2026-02-21 07:18:57"
$ws.Range("J18").Value = "1029.6 hPa"
$ws.Range("O18").Value = "0.6 °C"
$ws.Range("E19").Value = "2026-02-21 07:18:59"
$ws.Range("H19").Value = "89%"
$ws.Range("K19").Value = "0.0 MJ/m2"
$ws.Range("E20").Value = "2026-02-21 07:19:02"
$ws.Range("H20").Value = "46%"
$ws.Range("M20").Value = "2.5 °C 6:53 TU"
$ws.Range("O20").Value = "0.0 °C"
$ws.Range("E21").Value = "2026-02-21 07:19:04"
$ws.Range("N21").Value = "-1.0 °C 6:43 TU"
$ws.Range("E22").Value = "2026-02-21 07:19:06"
$ws.Range("O22").Value = "-0.8 °C"
$ws.Range("E23").Value = "2026-02-21 07:19:09"
$ws.Range("H23").Value = "34%"
$ws.Range("M23").Value = "1.9 °C 6:47 TU"
$ws.Range("O23").Value = "0.4 °C"
$ws.Range("E24").Value = "2026-02-21 07:19:11"
$ws.Range("J24").Value = "1031.5 hPa"
$ws.Range("O24").Value = "1.1 °C"
$ws.Range("E25").Value = "2026-02-21 07:19:14"
$ws.Range("H25").Value = "39%"
$ws.Range("K25").Value = "0.0 MJ/m2"
$ws.Range("O25").Value = "0.4 °C"
$ws.Range("E26").Value = "2026-02-21 07:19:16"
$ws.Range("J26").Value = "1027.4 hPa"
$ws.Range("N26").Value = "5.2 °C 6:49 TU"
$ws.Range("O26").Value = "6.7 °C"
$ws.Range("E27").Value = "2026-02-21 07:19:18"
$ws.Range("H27").Value = "39%"
$ws.Range("M27").Value = "3.6 °C 6:56 TU"
$ws.Range("O27").Value = "1.1 °C"
$ws.Range("E28").Value = "2026-02-21 07:19:21"
$ws.Range("J28").Value = "1030.5 hPa"
$ws.Range("N28").Value = "-1.1 °C 6:54 TU"
$ws.Range("O28").Value = "0.3 °C"
$ws.Range("E29").Value = "2026-02-21 07:19:23"
$ws.Range("H29").Value = "68%"
$ws.Range("N29").Value = "7.5 °C 6:38 TU"
$ws.Range("O29").Value = "9.4 °C"
$ws.Range("E30").Value = "2026-02-21 07:19:25"
$ws.Range("H30").Value = "74%"
$ws.Range("J30").Value = "1028.6 hPa"
$ws.Range("K30").Value = "0.0 MJ/m2"
$ws.Range("N30").Value = "4.4 °C 6:38 TU"
$ws.Range("O30").Value = "8.4 °C"
$ws.Range("E31").Value = "2026-02-21 07:19:28"
$ws.Range("H31").Value = "59%"
$ws.Range("J31").Value = "1027.1 hPa"
$ws.Range("K31").Value = "0.0 MJ/m2"
$ws.Range("M31").Value = "11.3 °C 6:51 TU"
$ws.Range("O31").Value = "10.4 °C"
$ws.Range("E32").Value = "2026-02-21 07:19:30"
$ws.Range("N32").Value = "-3.2 °C 6:57 TU"
$ws.Range("O32").Value = "0.1 °C"
$ws.Range("E33").Value = "2026-02-21 07:19:33"
$ws.Range("H33").Value = "79%"
$ws.Range("J33").Value = "1033.8 hPa"
$ws.Range("E34").Value = "2026-02-21 07:19:35"
$ws.Range("E35").Value = "2026-02-21 07:19:38"
$ws.Range("J35").Value = "1032.1 hPa"
$ws.Range("N35").Value = "2.4 °C 6:48 TU"
$ws.Range("E36").Value = "2026-02-21 07:19:40"
$ws.Range("H36").Value = "55%"
$ws.Range("J36").Value = "1028.2 hPa"
$ws.Range("N36").Value = "7.9 °C 6:49 TU"
$ws.Range("O36").Value = "11.6 °C"
$ws.Range("E37").Value = "2026-02-21 07:19:43"
$ws.Range("J37").Value = "1033.1 hPa"
$ws.Range("E38").Value = "2026-02-21 07:19:45"
$ws.Range("N38").Value = "1.9 °C 6:59 TU"
$ws.Range("O38").Value = "3.9 °C"
$ws.Range("E39").Value = "2026-02-21 07:19:47"
$ws.Range("K39").Value = "0.0 MJ/m2"
$ws.Range("N39").Value = "-0.6 °C 6:34 TU"
$ws.Range("E40").Value = "2026-02-21 07:19:50"
$ws.Range("H40").Value = "68%"
$ws.Range("O40").Value = "3.1 °C"
$ws.Range("E41").Value = "2026-02-21 07:19:52"
$ws.Range("H41").Value = "76%"
$ws.Range("J41").Value = "1028.8 hPa"
$ws.Range("N41").Value = "3.7 °C 6:35 TU"
$ws.Range("O41").Value = "7.3 °C"
$ws.Range("E42").Value = "2026-02-21 07:19:55"
$ws.Range("H42").Value = "81%"
$ws.Range("O42").Value = "7.1 °C"
$ws.Range("E43").Value = "2026-02-21 07:19:57"
$ws.Range("O43").Value = "0.0 °C"
$ws.Range("E44").Value = "2026-02-21 07:19:59"
$ws.Range("E45").Value = "2026-02-21 07:20:02"
$ws.Range("J45").Value = "1035.3 hPa"
$ws.Range("O45").Value = "0.1 °C"
$ws.Range("E46").Value = "2026-02-21 07:20:04"
$ws.Range("H46").Value = "87%"
$ws.Range("J46").Value = "1031.5 hPa"
$ws.Range("N46").Value = "1.7 °C 6:42 TU"
$ws.Range("O46").Value = "4.7 °C"
